$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of three species-observation records (rows 23-25):
#   new row 23 <- old row 25 data
#   new row 24 <- old row 23 data
#   new row 25 <- old row 24 data
# Only the columns that actually differ between the three original rows are
# touched (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
# Auktor, Aktivitet, Ost, Nord, Noggrannhet, Observatorer).

# Column A - Id
$ws.Range("A23").Value = 111941304
$ws.Range("A24").Value = 111941831
$ws.Range("A25").Value = 111941827

# Column B - Taxonsorteringsordning
$ws.Range("B23").Value = 78081
$ws.Range("B24").Value = 56398
$ws.Range("B25").Value = 77268

# Column E - TaxonId
$ws.Range("E23").Value = 229821
$ws.Range("E24").Value = 100109
$ws.Range("E25").Value = 228912

# Column F - Artnamn
$ws.Range("F23").Value = "Vedflamlav"
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("F25").Value = "Mörk kolflarnlav"

# Column G - Vetenskapligt namn
$ws.Range("G23").Value = "Ramboldia elabens"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("G25").Value = "Carbonicola myrmecina"

# Column H - Auktor
$ws.Range("H23").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("H25").Value = "(Ach.) Bendiksby & Timdal"

# Column M - Aktivitet (only old row 23 had a value here: "aldre spar")
$ws.Range("M23").ClearContents()
$ws.Range("M24").Value = "äldre spår"
$ws.Range("M25").ClearContents()

# Column Q - Ost
$ws.Range("Q23").Value = 466297.5338563451
$ws.Range("Q24").Value = 466322.7466770636
$ws.Range("Q25").Value = 466325.0201382869

# Column R - Nord
$ws.Range("R23").Value = 6820498.775792331
$ws.Range("R24").Value = 6821027.758616986
$ws.Range("R25").Value = 6821015.807131591

# Column S - Noggrannhet
$ws.Range("S23").Value = 25
$ws.Range("S24").Value = 10
$ws.Range("S25").Value = 10

# Column AX - Observatorer
$ws.Range("AX23").Value = "Bengt Oldhammer, Birgitta Kvist, Peter Turander"
$ws.Range("AX24").Value = "Bengt Oldhammer, Peter Turander, Birgitta Kvist"
$ws.Range("AX25").Value = "Bengt Oldhammer, Peter Turander, Birgitta Kvist"
